$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the empty "experimentObservations" column (E), shifting F:J left to E:I
$ws.Range("E1").EntireColumn.Delete()

# Correct the strain value from KN99alpha to TDY451 (now in column E after the shift)
$ws.Range("E2:E7").Value = "TDY451"

# Update the selection to mirror the saved view state
$ws.Range("E1:E1048576").Select()
